$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.748.84'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '3.447.10'
$ws.Range("E3").Value = '  +2.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.44'
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.70'
$ws.Range("E6").Value = '  +8.03%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.447.73'
$ws.Range("E7").Value = '  +2.24%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.67'
$ws.Range("E10").Value = '  +3.34%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").Value = '4.033.72'
$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("E14").Value = '  -1.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.26'
$ws.Range("E15").Value = '  +5.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").Value = '  +0.43%  '

$ws.Range("D17").Value = '3.450.27'
$ws.Range("E17").Value = '  +2.31%  '

$ws.Range("D18").Value = '61.827.59'
$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.12'
$ws.Range("E19").Value = '  +3.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.54'
$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.59'
$ws.Range("E22").Value = '  +1.91%  '

$ws.Range("E23").Value = '  +1.16%  '

$ws.Range("D24").Value = '3.570.28'
$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.02'
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '72.01'
$ws.Range("E26").Value = '  +0.84%  '

$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000126'
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("E28").Value = '  +10.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.76'
$ws.Range("E29").Value = '  +2.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.61'
$ws.Range("E30").Value = '  -8.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("E33").Value = '  +0.92%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.88'
$ws.Range("E35").Value = '  +1.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.30'
$ws.Range("E36").Value = '  +1.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.02'
$ws.Range("E37").Value = '  +2.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.57'
$ws.Range("E38").Value = '  +2.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '166.63'
$ws.Range("E39").Value = '  +1.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0793'
$ws.Range("E40").Value = '  +2.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.34'
$ws.Range("E41").Value = '  +7.32%  '

$ws.Range("E42").Value = '  +2.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.74'
$ws.Range("E43").Value = '  +1.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.20'
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("E46").Value = '  +1.98%  '

$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("D48").Value = '2.632.95'
$ws.Range("E48").Value = '  +12.00%  '

$ws.Range("E49").Value = '  +5.52%  '

$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.20'
$ws.Range("E51").Value = '  +8.47%  '
